$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix 0x201 Byte Order: the "Byte #" header row should count up from 0 (LSB, rightmost byte)
# to 7 (MSB, leftmost byte) left-to-right, instead of counting down from 7 to 0.
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
$ws.Range("G1").Value = 5
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 7

# Restore the cursor/selection position as left by the author
[void]$ws.Range("B4").Select()
